$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D, shifting the existing quarterly data
# (columns D:K) two columns to the right (to F:M)
$ws.Range("D7:E7").EntireColumn.Insert()

# Copy number formatting from column F (the former column D) into the new D:E
# columns, for each contiguous block of data rows
$ws.Range("F7:F35").Copy()
$ws.Range("D7:E35").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Range("F38:F77").Copy()
$ws.Range("D38:E77").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Range("F80:F102").Copy()
$ws.Range("D80:E102").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# Populate the two new quarter columns (D = period ending 2019-01-31,
# E = period ending 2018-10-31) with data
$ws.Range("D7").Value = 43496
$ws.Range("E7").Value = 43404
$ws.Range("D8").Value = 1290600
$ws.Range("E8").Value = 1756000
$ws.Range("D9").Value = 1149000
$ws.Range("E9").Value = 1548700
$ws.Range("D10").Value = 141600
$ws.Range("E10").Value = 207300
$ws.Range("D12:E12").Value = "NA"
$ws.Range("D13:E13").Value = 0
$ws.Range("D14").Value = 42100
$ws.Range("E14").Value = 57100
$ws.Range("D15").Value = 12500
$ws.Range("E15").Value = 12600
$ws.Range("D17").Value = 1288700
$ws.Range("E17").Value = 1724500
$ws.Range("D18").Value = 1900
$ws.Range("E18").Value = 31500
$ws.Range("D20:E20").Value = 0
$ws.Range("D21").Value = 25400
$ws.Range("E21").Value = 54600
$ws.Range("D22:E22").Value = 0
$ws.Range("D23").Value = 1900
$ws.Range("E23").Value = 31500
$ws.Range("D24").Value = 7300
$ws.Range("E24").Value = 17600
$ws.Range("D25:E25").Value = 0
$ws.Range("D26").Value = -5400
$ws.Range("E26").Value = 14000
$ws.Range("D27").Value = -5400
$ws.Range("E27").Value = 14000
$ws.Range("D28:E28").Value = 0
$ws.Range("D29:E29").Value = "NA"
$ws.Range("D30:E30").Value = 0
$ws.Range("D31:E31").Value = 0
$ws.Range("D32:E32").Value = 0
$ws.Range("D33").Value = -5400
$ws.Range("E33").Value = 14000
$ws.Range("D34:E34").Value = 0
$ws.Range("D35").Value = -5400
$ws.Range("E35").Value = 14000
$ws.Range("D38").Value = 43496
$ws.Range("E38").Value = 43404
$ws.Range("D41").Value = 305800
$ws.Range("E41").Value = 224900
$ws.Range("D42:E42").Value = 0
$ws.Range("D43").Value = 344800
$ws.Range("E43").Value = 503800
$ws.Range("D44").Value = 561800
$ws.Range("E44").Value = 565300
$ws.Range("D45").Value = 35700
$ws.Range("E45").Value = 30900
$ws.Range("D46").Value = 1248200
$ws.Range("E46").Value = 1325000
$ws.Range("D47").Value = 48300
$ws.Range("E47").Value = 47000
$ws.Range("D48").Value = 550500
$ws.Range("E48").Value = 543700
$ws.Range("D49").Value = 740900
$ws.Range("E49").Value = 753500
$ws.Range("D50:E50").Value = 0
$ws.Range("D51:E51").Value = 0
$ws.Range("D52").Value = 142400
$ws.Range("E52").Value = 131300
$ws.Range("D53:E53").Value = 0
$ws.Range("D54").Value = 2730200
$ws.Range("E54").Value = 2800400
$ws.Range("D57").Value = 219900
$ws.Range("E57").Value = 255500
$ws.Range("D58:E58").Value = 0
$ws.Range("D59").Value = 530800
$ws.Range("E59").Value = 544300
$ws.Range("D60").Value = 750600
$ws.Range("E60").Value = 799800
$ws.Range("D61:E61").Value = 0
$ws.Range("D62").Value = 72600
$ws.Range("E62").Value = 72300
$ws.Range("D63:E63").Value = 0
$ws.Range("D64:E64").Value = 0
$ws.Range("D65:E65").Value = 0
$ws.Range("D66").Value = 823300
$ws.Range("E66").Value = 872100
$ws.Range("D68:E68").Value = 0
$ws.Range("D69:E69").Value = 0
$ws.Range("D70:E70").Value = 0
$ws.Range("D71:E71").Value = 0
$ws.Range("D72").Value = 1984900
$ws.Range("E72").Value = 2010900
$ws.Range("D73:E73").Value = 0
$ws.Range("D74:E74").Value = 0
$ws.Range("D75:E75").Value = 0
$ws.Range("D76").Value = 1906900
$ws.Range("E76").Value = 1928300
$ws.Range("D77:E77").Value = 0
$ws.Range("D80").Value = 43496
$ws.Range("E80").Value = 43404
$ws.Range("D81").Value = -5400
$ws.Range("E81").Value = 14000
$ws.Range("D83").Value = 23500
$ws.Range("E83").Value = 23100
$ws.Range("D84:E84").Value = 0
$ws.Range("D85:E85").Value = 0
$ws.Range("D86:E86").Value = 0
$ws.Range("D87:E87").Value = 0
$ws.Range("D88:E88").Value = 0
$ws.Range("D89").Value = 150500
$ws.Range("E89").Value = -15800
$ws.Range("D91").Value = -20300
$ws.Range("E91").Value = -34500
$ws.Range("D92:E92").Value = 0
$ws.Range("D93:E93").Value = 0
$ws.Range("D94").Value = -23800
$ws.Range("E94").Value = -34400
$ws.Range("D96:E96").Value = 0
$ws.Range("D97:E97").Value = 0
$ws.Range("D98:E98").Value = 0
$ws.Range("D99:E99").Value = 0
$ws.Range("D100").Value = -45700
$ws.Range("E100").Value = -100
$ws.Range("D101:E101").Value = 0
$ws.Range("D102").Value = 80900
$ws.Range("E102").Value = -50300
